$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newParticipantsQuery = "MATCH (p:participant)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nOPTIONAL MATCH (samp)<--(f:file)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nWITH s, p, samp, f, g, diag`nWHERE g.platform in ['Illumina HiSeq X Ten']`nwith p`nOPTIONAL MATCH (p)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nWITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`nRETURN`ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY p.participant_id LIMIT 100"

$ws.Range("B2").Value = $newParticipantsQuery

# Update the window view to match the committed state: scrolled down a bit,
# with B5 selected instead of C2.
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("B5").Select()
